$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so values like "585.53" are not
# auto-converted to numbers by Excel (matches original inlineStr/text cells).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '66.510.41'
$ws.Range('E2').Value = '  -0.94%  '
$ws.Range('D3').Value = '3.450.21'
$ws.Range('E3').Value = '  -0.39%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = '585.53'
$ws.Range('E5').Value = '  -1.39%  '
$ws.Range('D6').Value = '176.20'
$ws.Range('E6').Value = '  -2.26%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('D8').Value = '0.601'
$ws.Range('E8').Value = '  -1.27%  '
$ws.Range('D9').Value = '3.445.70'
$ws.Range('E9').Value = '  -0.50%  '
$ws.Range('E10').Value = '  -6.41%  '
$ws.Range('D11').Value = '6.87'
$ws.Range('E11').Value = '  -1.23%  '
$ws.Range('D12').Value = '0.416'
$ws.Range('E12').Value = '  -3.49%  '
$ws.Range('D13').Value = '4.045.40'
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('B14').Value = 'TRON'
$ws.Range('C14').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D14').Value = '0.133'
$ws.Range('E14').Value = '  -0.25%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').Value = '30.21'
$ws.Range('E15').Value = '  -5.60%  '
$ws.Range('D16').Value = '66.506.62'
$ws.Range('E16').Value = '  -0.96%  '
$ws.Range('E17').Value = '  -2.91%  '
$ws.Range('D18').Value = '3.450.76'
$ws.Range('E18').Value = '  -0.18%  '
$ws.Range('D19').Value = '5.95'
$ws.Range('E19').Value = '  -4.69%  '
$ws.Range('D20').Value = '13.80'
$ws.Range('E20').Value = '  -2.94%  '
$ws.Range('D21').Value = '378.31'
$ws.Range('E21').Value = '  -3.22%  '
$ws.Range('E22').Value = '  -1.30%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('D24').Value = '72.36'
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('D25').Value = '5.73'
$ws.Range('E25').Value = '  -0.38%  '
$ws.Range('D26').Value = '0.536'
$ws.Range('E26').Value = '  -0.26%  '
$ws.Range('E27').Value = '  -1.62%  '
$ws.Range('D28').Value = '9.77'
$ws.Range('E28').Value = '  -5.90%  '
$ws.Range('E29').Value = '  +0.72%  '
$ws.Range('E30').Value = '  +0.08%  '
$ws.Range('D31').Value = '24.24'
$ws.Range('E31').Value = '  +3.29%  '
$ws.Range('D32').Value = '5.84'
$ws.Range('E32').Value = '  -4.89%  '
$ws.Range('E33').Value = '  -2.85%  '
$ws.Range('B34').Value = 'USDe'
$ws.Range('C34').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.06%  '
$ws.Range('B35').Value = 'Fetch.AI'
$ws.Range('C35').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D35').Value = '1.31'
$ws.Range('E35').Value = '  -6.45%  '
$ws.Range('D36').Value = '7.17'
$ws.Range('E36').Value = '  -2.29%  '
$ws.Range('D37').Value = '1.57'
$ws.Range('E37').Value = '  -0.79%  '
$ws.Range('D38').Value = '161.30'
$ws.Range('E38').Value = '  -1.44%  '
$ws.Range('D39').Value = '29.34'
$ws.Range('E39').Value = '  +12.30%  '
$ws.Range('D40').Value = '0.890'
$ws.Range('E40').Value = '  +1.72%  '
$ws.Range('B41').Value = 'dogwifhat'
$ws.Range('C41').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D41').Value = '2.63'
$ws.Range('E41').Value = '  -6.80%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '1.78'
$ws.Range('E42').Value = '  -4.75%  '
$ws.Range('D43').Value = '4.50'
$ws.Range('E43').Value = '  -3.58%  '
$ws.Range('D44').Value = '2.736.04'
$ws.Range('E44').Value = '  -0.24%  '
$ws.Range('D45').Value = '6.36'
$ws.Range('E45').Value = '  -6.54%  '
$ws.Range('D46').Value = '0.0695'
$ws.Range('E46').Value = '  -3.55%  '
$ws.Range('D47').Value = '40.68'
$ws.Range('E47').Value = '  -1.60%  '
$ws.Range('D48').Value = '24.52'
$ws.Range('E48').Value = '  -6.65%  '
$ws.Range('E49').Value = '  -1.71%  '
$ws.Range('D50').Value = '309.13'
$ws.Range('E50').Value = '  -5.63%  '
$ws.Range('D51').Value = '0.827'
$ws.Range('E51').Value = '  -1.59%  '

# Restore default General format / style on column D so no stray styling
# is introduced (matches the unstyled text cells in the source workbook).
$ws.Range("D2:D51").NumberFormat = "General"
$ws.Range("D2:D51").Style = "Normal"
